$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1496.125
$ws.Range("I15").Value = 1496.125
$ws.Range("K15").Value = 4488.375
$ws.Range("M15").Value = -4319.375
$ws.Range("H69").Value = 23815614
$ws.Range("I69").Value = 83336830
$ws.Range("J69").Value = 7125
$ws.Range("K69").Value = 250010490
$ws.Range("L69").Value = 21375
$ws.Range("M69").Value = -250009616
$ws.Range("N69").Value = -23123
$ws.Range("H72").Value = 23815614
$ws.Range("I72").Value = 83336830
$ws.Range("J72").Value = 7125
$ws.Range("K72").Value = 750031470
$ws.Range("L72").Value = 64125
$ws.Range("M72").Value = -750027102
$ws.Range("N72").Value = -72861
$ws.Range("H82").Value = 6766.9375
$ws.Range("I82").Value = 1975.2858
$ws.Range("K82").Value = 5925.857400000001
$ws.Range("M82").Value = -5519.857400000001
$ws.Range("H85").Value = 6766.9375
$ws.Range("I85").Value = 1975.2858
$ws.Range("K85").Value = 5925.857400000001
$ws.Range("M85").Value = -4521.857400000001
$ws.Range("H86").Value = 4461.125
$ws.Range("I86").Value = 4098.4287
$ws.Range("K86").Value = 4098.4287
$ws.Range("M86").Value = -2975.4287
$ws.Range("H89").Value = 4461.125
$ws.Range("I89").Value = 4098.4287
$ws.Range("K89").Value = 20492.1435
$ws.Range("M89").Value = -14876.1435
$ws.Range("H112").Value = 3298.4473
$ws.Range("J112").Value = 3360.5945
$ws.Range("L112").Value = 10081.7835
$ws.Range("N112").Value = -12297.7835
$ws.Range("H121").Value = 4966.2383
$ws.Range("J121").Value = 4966.2383
$ws.Range("L121").Value = 14898.7149
$ws.Range("N121").Value = -18392.7149
$ws.Range("H132").Value = 4007.3513
$ws.Range("I132").Value = 2071.3438
$ws.Range("K132").Value = 6214.0314
$ws.Range("M132").Value = -3684.0314
$ws.Range("H137").Value = 42921420
$ws.Range("J137").Value = 83339170
$ws.Range("L137").Value = 250017510
$ws.Range("N137").Value = -250022610
$ws.Range("H138").Value = 2758.06
$ws.Range("I138").Value = 1145.96
$ws.Range("J138").Value = 3295.4268
$ws.Range("K138").Value = 3437.88
$ws.Range("L138").Value = 9886.2804
$ws.Range("M138").Value = 1702.12
$ws.Range("N138").Value = -20166.2804

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2443.7188
$ws.Range("I45").Value = 1912.24
$ws.Range("K45").Value = 1912.24
$ws.Range("M45").Value = -1535.24
$ws.Range("H74").Value = 7813371
$ws.Range("I74").Value = 9616236
$ws.Range("J74").Value = 956.3333
$ws.Range("K74").Value = 9616236
$ws.Range("L74").Value = 956.3333
$ws.Range("M74").Value = -9615362
$ws.Range("N74").Value = -2704.3333
$ws.Range("H77").Value = 7813371
$ws.Range("I77").Value = 9616236
$ws.Range("J77").Value = 956.3333
$ws.Range("K77").Value = 48081180
$ws.Range("L77").Value = 4781.6665
$ws.Range("M77").Value = -48076812
$ws.Range("N77").Value = -13517.6665

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 14660
$ws.Range("I5").Value = 14660
$ws.Range("K5").Value = 14660
$ws.Range("M5").Value = -14547
$ws.Range("H86").Value = 3520.2083
$ws.Range("I86").Value = 3567.4167
$ws.Range("J86").Value = 3473
$ws.Range("K86").Value = 3567.4167
$ws.Range("L86").Value = 3473
$ws.Range("M86").Value = -2444.4167
$ws.Range("N86").Value = -5719
$ws.Range("H89").Value = 3520.2083
$ws.Range("I89").Value = 3567.4167
$ws.Range("J89").Value = 3473
$ws.Range("K89").Value = 17837.0835
$ws.Range("L89").Value = 17365
$ws.Range("M89").Value = -12221.0835
$ws.Range("N89").Value = -28597

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 129.6
$ws.Range("I7").Value = 172.125
$ws.Range("K7").Value = 172.125
$ws.Range("M7").Value = -59.125
$ws.Range("H76").Value = 6333
$ws.Range("I76").Value = 6333
$ws.Range("K76").Value = 6333
$ws.Range("M76").Value = -6018
$ws.Range("H79").Value = 6333
$ws.Range("I79").Value = 6333
$ws.Range("K79").Value = 6333
$ws.Range("M79").Value = -5241
$ws.Range("H99").Value = 6271.2085
$ws.Range("I99").Value = 6735.88
$ws.Range("J99").Value = 5766.1304
$ws.Range("K99").Value = 6735.88
$ws.Range("L99").Value = 5766.1304
$ws.Range("M99").Value = -5237.88
$ws.Range("N99").Value = -8762.1304
$ws.Range("H126").Value = 6271.2085
$ws.Range("I126").Value = 6735.88
$ws.Range("J126").Value = 5766.1304
$ws.Range("K126").Value = 20207.64
$ws.Range("L126").Value = 17298.3912
$ws.Range("M126").Value = -17737.64
$ws.Range("N126").Value = -22238.3912
$ws.Range("H132").Value = 23257618
$ws.Range("I132").Value = 28572870
$ws.Range("K132").Value = 85718610
$ws.Range("M132").Value = -85716080

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2076.6
$ws.Range("I68").Value = 1018.5
$ws.Range("K68").Value = 3055.5
$ws.Range("M68").Value = -2244.5
$ws.Range("H71").Value = 2076.6
$ws.Range("I71").Value = 1018.5
$ws.Range("K71").Value = 9166.5
$ws.Range("M71").Value = -5110.5
$ws.Range("H75").Value = 1670.7142
$ws.Range("I75").Value = 449.5
$ws.Range("J75").Value = 2159.2
$ws.Range("K75").Value = 1348.5
$ws.Range("L75").Value = 6477.599999999999
$ws.Range("M75").Value = -350.5
$ws.Range("N75").Value = -8473.599999999999
$ws.Range("H78").Value = 1670.7142
$ws.Range("I78").Value = 449.5
$ws.Range("J78").Value = 2159.2
$ws.Range("K78").Value = 4045.5
$ws.Range("L78").Value = 19432.8
$ws.Range("M78").Value = 946.5
$ws.Range("N78").Value = -29416.8
$ws.Range("H137").Value = 75002570
$ws.Range("J137").Value = 4021.6667
$ws.Range("L137").Value = 12065.0001
$ws.Range("N137").Value = -22265.0001
$ws.Range("H138").Value = 3105.9583
$ws.Range("I138").Value = 3130.8096
$ws.Range("K138").Value = 9392.4288
$ws.Range("M138").Value = -4252.4288

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5845.154
$ws.Range("I70").Value = 5082.1665
$ws.Range("K70").Value = 5082.1665
$ws.Range("M70").Value = -4812.1665
$ws.Range("H73").Value = 5845.154
$ws.Range("I73").Value = 5082.1665
$ws.Range("K73").Value = 5082.1665
$ws.Range("M73").Value = -4146.1665
$ws.Range("H102").Value = 9627559
$ws.Range("I102").Value = 11595469
$ws.Range("K102").Value = 11595469
$ws.Range("M102").Value = -11593847
$ws.Range("H122").Value = 298139.94
$ws.Range("I122").Value = 557628.7
$ws.Range("K122").Value = 1672886.1
$ws.Range("M122").Value = -1670436.1
$ws.Range("H123").Value = 44717.57
$ws.Range("J123").Value = 44717.57
$ws.Range("L123").Value = 44717.57
$ws.Range("N123").Value = -49617.57
$ws.Range("H132").Value = 93783.77
$ws.Range("I132").Value = 143808.86
$ws.Range("K132").Value = 431426.58
$ws.Range("M132").Value = -428896.58

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 20837474
$ws.Range("I40").Value = 3026.889
$ws.Range("K40").Value = 3026.889
$ws.Range("M40").Value = -2890.889
$ws.Range("H82").Value = 2022.4482
$ws.Range("I82").Value = 1496.3334
$ws.Range("J82").Value = 2883.3635
$ws.Range("K82").Value = 1496.3334
$ws.Range("L82").Value = 2883.3635
$ws.Range("M82").Value = -1135.3334
$ws.Range("N82").Value = -3605.3635
$ws.Range("H85").Value = 2022.4482
$ws.Range("I85").Value = 1496.3334
$ws.Range("J85").Value = 2883.3635
$ws.Range("K85").Value = 1496.3334
$ws.Range("L85").Value = 2883.3635
$ws.Range("M85").Value = -248.3334
$ws.Range("N85").Value = -5379.363499999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3089.0571
$ws.Range("I136").Value = 1754.9259
$ws.Range("K136").Value = 5264.7777
$ws.Range("M136").Value = -2714.7777
$ws.Range("H140").Value = 99997.336
$ws.Range("J140").Value = 99997.336
$ws.Range("L140").Value = 99997.336
$ws.Range("N140").Value = -110357.336
